$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: averages under the PC Video Process Time table ---
$ws.Range("B9").Formula = "=AVERAGE(B6:B8)"
$ws.Range("C9").Formula = "=AVERAGE(C6:C8)"
$ws.Range("D9").Formula = "=AVERAGE(D6:D8)"

# --- Fix up the old "merge1/2/3" labels and add a new one ---
$ws.Range("A14").Value = "Merge"
$ws.Range("A15").Value = "3worker"
$ws.Range("A16").Value = "8worker"
$ws.Range("A17").Value = "10worker"

# --- New "Process Time" section (rows 19-26) ---
$ws.Range("A19").Value = "Process Time"

$ws.Range("A20").HorizontalAlignment = -4131
$ws.Range("B20").Value = "3Worker"
$ws.Range("C20").Value = "8Worker"
$ws.Range("D20").Value = "12Worker"
$ws.Range("E20").Value = "Server"
$ws.Range("I20").Value = "MergeSpeed"
$ws.Range("J20").Value = "4 to 2"
$ws.Range("K20").Value = "2 to 1"
$ws.Range("J20:O20").HorizontalAlignment = -4108

$ws.Range("A21").Value = "360 Convert"
$ws.Range("B21").Value = "1.61,2.71;1.38,2.16"
$ws.Range("C21").Value = "1.61,2.71;1.38,2.17"
$ws.Range("D21").Value = "1.61,2.71;1.38,2.18"
$ws.Range("E21").Value = 5.99
$ws.Range("I21").Value = 360

$ws.Range("A22").Value = "360 Merge"
$ws.Range("B22").Value = 7.8
$ws.Range("I22").Value = 1080

$ws.Range("A23").Value = "1080 Convert"
$ws.Range("B23").Value = "32.71,23.16;22.64,30.28;22.90,30.68"
$ws.Range("C23").Value = "32.71,23.16;22.64,30.29"
$ws.Range("D23").Value = "32.71,23.16;22.64,30.30"
$ws.Range("E23").Value = "40.9,43.37"
$ws.Range("I23").Value = "4k"

$ws.Range("A24").Value = "1080 Merge"

$ws.Range("A25").Value = "4k Convert"
$ws.Range("B25").Value = "141.12,133.12;120.22,120.53"
$ws.Range("C25").Value = "141.12,133.12;120.22,120.54"
$ws.Range("D25").Value = "141.12,133.12;120.22,120.55"
$ws.Range("E25").Value = "194.22,183.26"

$ws.Range("A26").Value = "4K Merge"

$ws.Range("E21:E22").HorizontalAlignment = -4108
$ws.Range("E21:E22").Merge()
$ws.Range("E23:E24").HorizontalAlignment = -4108
$ws.Range("E23:E24").Merge()
$ws.Range("E25:E26").HorizontalAlignment = -4108
$ws.Range("E25:E26").Merge()

$ws.Range("I21:I23").HorizontalAlignment = -4108
$ws.Range("I21:I23").VerticalAlignment = -4108

# --- "360/1080/4K Tramission Time" tables (rows 29-45) ---
$ws.Range("A29").Value = "360 Tramission Time"
$ws.Range("J29:O29").HorizontalAlignment = -4108
$ws.Range("J29:K29").Merge()
$ws.Range("L29:M29").Merge()
$ws.Range("N29:O29").Merge()

$ws.Range("A30").Value = "Worker"
$ws.Range("B30").Value = "3Worker"
$ws.Range("C30").Value = "8Worker"
$ws.Range("D30").Value = "12Worker"
$ws.Range("E30").Value = "Server"

$ws.Range("A31").Value = "2Requester"
$ws.Range("E31").Formula = "=3.04+3.04"

$ws.Range("A32").Value = "8Requester"
$ws.Range("E32").Formula = "=6.1"

$ws.Range("A33").Value = "12Requertser"
$ws.Range("E33").Formula = "=3.02+12*3.02/5.12"

$ws.Range("A35").Value = "1080 Tramission Time"
$ws.Range("A36").Value = "Worker"
$ws.Range("B36").Value = "3Worker"
$ws.Range("C36").Value = "8Worker"
$ws.Range("D36").Value = "12Worker"
$ws.Range("E36").Value = "Server"
$ws.Range("A37").Value = "2Requester"
$ws.Range("A38").Value = "8Requester"
$ws.Range("A39").Value = "12Requertser"

$ws.Range("A41").Value = "4K Tramission Time"
$ws.Range("A42").Value = "Worker"
$ws.Range("B42").Value = "3Worker"
$ws.Range("C42").Value = "8Worker"
$ws.Range("D42").Value = "12Worker"
$ws.Range("E42").Value = "Server"
$ws.Range("A43").Value = "2Requester"
$ws.Range("A44").Value = "8Requester"
$ws.Range("A45").Value = "12Requertser"

# --- size / new_* block (rows 48-51) ---
$ws.Range("B48").Value = "size"
$ws.Range("A49").Value = "new_360"
$ws.Range("B49").Value = 2.65
$ws.Range("A50").Value = "new_1080"
$ws.Range("B50").Value = 68.6
$ws.Range("A51").Value = "new_4k"
$ws.Range("B51").Value = 149

$ws.Columns.Item(9).ColumnWidth = 16.6640625

$ws.Range("J23").Select()
